# Update RF012 - Ajuda workbook: swap the order of the "Niveis das
# Competencias" and "Avaliacoes" test-case blocks (TC1 becomes Avaliacoes,
# TC2 becomes Niveis das Competencias), and fix wording from "cadastradas"
# to "cadastrados" for the Niveis das Competencias listing text.
# From 1.3 to 1.4 version

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block 1 (rows 10-12, under "TC1") now describes "Avaliacoes" ---
$ws.Range("B10").Value = "Lider de Pessoas acessa a funcionalidade de gestao de Avaliacoes a partir do menu inicial"
$ws.Range("D10").Value = "SYSTEM exibe a listagem das Avaliacoes cadastradas com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
$ws.Range("B11").Value = "Lider de Pessoas clica na opcao 'Ajuda' para visualizar a ajuda da gestao de Avaliacoes"
# D11, A12, B12 (shared help panel / postcondition text) stay unchanged

# --- Block 2 (rows 19-21, under "TC2") now describes "Niveis das Competencias" ---
$ws.Range("B19").Value = "Lider de Pessoas acessa a funcionalidade de gestao de Niveis das Competencias a partir do menu inicial"
$ws.Range("D19").Value = "SYSTEM exibe a listagem dos Niveis das Competencias cadastrados com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
$ws.Range("B20").Value = "Lider de Pessoas clica na opcao 'Ajuda' para visualizar a ajuda da gestao de Niveis das Competencias"
# D20, A21, B21 (shared help panel / postcondition text) stay unchanged
